# "Generate Report for Archive"
#
# 1. Update the localization status text from "Ready for handoff" to
#    "In Translation" everywhere it appears (Overview!E2:F2, and the
#    "Status" column - column C - on the per-locale sheets).
# 2. Narrow the now-shorter "Status"-related columns (Overview columns
#    E/F, and column C on the zh-cn / de-de sheets) to match the new,
#    shorter header/value text.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = "In Translation"
$overview.Range("F2").Value = "In Translation"
$overview.Columns.Item(5).ColumnWidth = 12.5
$overview.Columns.Item(6).ColumnWidth = 12.5

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = "In Translation"
$zhcn.Columns.Item(3).ColumnWidth = 12.5

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = "In Translation"
$dede.Columns.Item(3).ColumnWidth = 12.5
